$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking account number as real text (keeping
# leading zeros) without leaving a lingering custom number-format style
# on the cell - matches how the original file stores these values
# (plain inlineStr cells with no style).
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Helper: find the worksheet row number whose column-A value equals a
# given account number (searches the data rows only).
function Find-RowByAccount($account) {
    $lastRow = $ws.UsedRange.Rows.Count
    for ($r = 2; $r -le $lastRow; $r++) {
        if ($ws.Cells.Item($r, 1).Value2 -eq $account) {
            return $r
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# 1) Remove the LEVI (005206566) and RACHEL (004588677) rows entirely.
#    Delete the lower row first so the earlier row's index stays valid.
# ---------------------------------------------------------------------
$rowRachel = Find-RowByAccount "004588677"
$ws.Range("A$($rowRachel):C$($rowRachel)").EntireRow.Delete()

$rowLevi = Find-RowByAccount "005206566"
$ws.Range("A$($rowLevi):C$($rowLevi)").EntireRow.Delete()

# ---------------------------------------------------------------------
# 2) Move the THOMAS/BLOCO/PHYLIA/BLUEMETRIX block so it sits right
#    before the PATRICIA (005255637) row instead of right after it.
# ---------------------------------------------------------------------
$accounts = @("004224011", "004364200", "004690692", "001761119")
$rowFirst = Find-RowByAccount $accounts[0]
$rowLast = $rowFirst + $accounts.Length - 1

$blockVals = @()
for ($r = $rowFirst; $r -le $rowLast; $r++) {
    $blockVals += ,@($ws.Cells.Item($r,1).Value2, $ws.Cells.Item($r,2).Value2, $ws.Cells.Item($r,3).Value2)
}

$delFrom = $ws.Rows.Item($rowFirst)
$delTo = $ws.Rows.Item($rowLast)
$ws.Range($delFrom, $delTo).Delete()

$rowPatricia = Find-RowByAccount "005255637"

$insFrom = $ws.Rows.Item($rowPatricia)
$insTo = $ws.Rows.Item($rowPatricia + $accounts.Length - 1)
$ws.Range($insFrom, $insTo).Insert()

for ($i = 0; $i -lt $accounts.Length; $i++) {
    $r = $rowPatricia + $i
    Set-TextValue $ws.Cells.Item($r,1) $blockVals[$i][0]
    $ws.Cells.Item($r,2).Value = $blockVals[$i][1]
    $ws.Cells.Item($r,3).Value = $blockVals[$i][2]
}

# ---------------------------------------------------------------------
# 3) PATRICIA's balance changes from 20000 to 4000 (she moved down by
#    the 4 rows inserted above).
# ---------------------------------------------------------------------
$rowPatricia = Find-RowByAccount "005255637"
$ws.Cells.Item($rowPatricia, 3).Value = 4000

# ---------------------------------------------------------------------
# 4) Add a new row for SERGIO (004975924 / 850.52) right after RODRIGO
#    (004392159).
# ---------------------------------------------------------------------
$rowRodrigo = Find-RowByAccount "004392159"
$rowSergio = $rowRodrigo + 1
$ws.Rows.Item($rowSergio).Insert()
Set-TextValue $ws.Cells.Item($rowSergio, 1) "004975924"
$ws.Cells.Item($rowSergio, 2).Value = "SERGIO"
$ws.Cells.Item($rowSergio, 3).Value = 850.52
